$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1159.5
$ws.Range("I15").Value = 1159.5
$ws.Range("K15").Value = 3478.5
$ws.Range("M15").Value = -3309.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 1475
$ws.Range("I69").Value = 950
$ws.Range("K69").Value = 2850
$ws.Range("M69").Value = -1976

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1400
$ws.Range("I70").Value = 1400
$ws.Range("K70").Value = 4200
$ws.Range("M70").Value = -3930

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 1475
$ws.Range("I72").Value = 950
$ws.Range("K72").Value = 8550
$ws.Range("M72").Value = -4182

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1400
$ws.Range("I73").Value = 1400
$ws.Range("K73").Value = 4200
$ws.Range("M73").Value = -3264

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3329.6667
$ws.Range("I106").Value = 3329.6667
$ws.Range("K106").Value = 3329.6667
$ws.Range("M106").Value = -2698.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1864.6666
$ws.Range("I107").Value = 1864.6666
$ws.Range("K107").Value = 1864.6666
$ws.Range("M107").Value = 55.33339999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4473.5
$ws.Range("J61").Value = 2841.6667
$ws.Range("L61").Value = 2841.6667
$ws.Range("N61").Value = -3265.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 854.5
$ws.Range("I74").Value = 854.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 854.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 19.5
$ws.Range("N74").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 854.5
$ws.Range("I77").Value = 854.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 4272.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 95.5
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 9867.583000000001
$ws.Range("I122").Value = 11495.223
$ws.Range("J122").Value = 4984.6665
$ws.Range("K122").Value = 34485.669
$ws.Range("L122").Value = 14953.9995
$ws.Range("M122").Value = -32035.669
$ws.Range("N122").Value = -19853.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2313.2666
$ws.Range("I132").Value = 1641.8334
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 4925.5002
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -2395.5002
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4473.5
$ws.Range("J136").Value = 2841.6667
$ws.Range("L136").Value = 8525.000100000001
$ws.Range("N136").Value = -13625.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2977.8333
$ws.Range("I20").Value = 2914.6
$ws.Range("J20").Value = 3294
$ws.Range("K20").Value = 2914.6
$ws.Range("L20").Value = 3294
$ws.Range("M20").Value = -2667.6
$ws.Range("N20").Value = -3788

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 100.333336
$ws.Range("J80").Value = 121.8
$ws.Range("L80").Value = 121.8
$ws.Range("N80").Value = -2117.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 100.333336
$ws.Range("J83").Value = 121.8
$ws.Range("L83").Value = 609
$ws.Range("N83").Value = -10593

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4622.273
$ws.Range("I107").Value = 4336.7144
$ws.Range("J107").Value = 5122
$ws.Range("K107").Value = 4336.7144
$ws.Range("L107").Value = 5122
$ws.Range("M107").Value = -2416.7144
$ws.Range("N107").Value = -8962

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5313
$ws.Range("I134").Value = 5394.3
$ws.Range("K134").Value = 16182.9
$ws.Range("M134").Value = -13647.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 441.8
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 441.8
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 441.8
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -667.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4374.353
$ws.Range("I31").Value = 1738.3
$ws.Range("J31").Value = 8140.143
$ws.Range("K31").Value = 1738.3
$ws.Range("L31").Value = 8140.143
$ws.Range("M31").Value = -1443.3
$ws.Range("N31").Value = -8730.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4374.353
$ws.Range("I34").Value = 1738.3
$ws.Range("J34").Value = 8140.143
$ws.Range("K34").Value = 1738.3
$ws.Range("L34").Value = 8140.143
$ws.Range("M34").Value = -1536.3
$ws.Range("N34").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2835.5
$ws.Range("I107").Value = 1494
$ws.Range("K107").Value = 1494
$ws.Range("M107").Value = 426

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -54920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4575.4614
$ws.Range("I132").Value = 3200.2
$ws.Range("J132").Value = 5435
$ws.Range("K132").Value = 9600.599999999999
$ws.Range("L132").Value = 16305
$ws.Range("M132").Value = -7070.599999999999
$ws.Range("N132").Value = -21365

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5597.2
$ws.Range("I134").Value = 5622.75
$ws.Range("J134").Value = 5495
$ws.Range("K134").Value = 16868.25
$ws.Range("L134").Value = 16485
$ws.Range("M134").Value = -14333.25
$ws.Range("N134").Value = -21555

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 91.666664
$ws.Range("J2").Value = 112.5
$ws.Range("L2").Value = 675
$ws.Range("N2").Value = -901

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 895
$ws.Range("J22").Value = 1782
$ws.Range("L22").Value = 5346
$ws.Range("N22").Value = -5684

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 895
$ws.Range("J27").Value = 1782
$ws.Range("L27").Value = 5346
$ws.Range("N27").Value = -5550

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 571.6667
$ws.Range("I50").Value = 497.5
$ws.Range("K50").Value = 1492.5
$ws.Range("M50").Value = -1011.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 571.6667
$ws.Range("I53").Value = 497.5
$ws.Range("K53").Value = 1492.5
$ws.Range("M53").Value = -1011.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 4825
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 939.5
$ws.Range("J122").Value = 1149.5
$ws.Range("L122").Value = 10345.5
$ws.Range("N122").Value = -15245.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8416.5
$ws.Range("I80").Value = 8164.6665
$ws.Range("K80").Value = 8164.6665
$ws.Range("M80").Value = -7166.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 8416.5
$ws.Range("I83").Value = 8164.6665
$ws.Range("K83").Value = 40823.3325
$ws.Range("M83").Value = -35831.3325

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2827
$ws.Range("I102").Value = 2827
$ws.Range("K102").Value = 2827
$ws.Range("M102").Value = -1205

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 95000
$ws.Range("J130").Value = 95000
$ws.Range("L130").Value = 95000
$ws.Range("N130").Value = -105040

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2940.7273
$ws.Range("I16").Value = 2597.75
$ws.Range("K16").Value = 2597.75
$ws.Range("M16").Value = -2427.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2500
$ws.Range("I136").Value = 2666.8333
$ws.Range("J136").Value = 1999.5
$ws.Range("K136").Value = 8000.499899999999
$ws.Range("L136").Value = 5998.5
$ws.Range("M136").Value = -5450.499899999999
$ws.Range("N136").Value = -11098.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3949
$ws.Range("I122").Value = 3933.6667
$ws.Range("K122").Value = 11801.0001
$ws.Range("M122").Value = -9351.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 80000
$ws.Range("J125").Value = 80000
$ws.Range("L125").Value = 80000
$ws.Range("N125").Value = -89840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3039.1365
$ws.Range("I132").Value = 2525.7144
$ws.Range("J132").Value = 3937.625
$ws.Range("K132").Value = 7577.1432
$ws.Range("L132").Value = 11812.875
$ws.Range("M132").Value = -5047.1432
$ws.Range("N132").Value = -16872.875
